$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "FAPs"
$ws.Cells.Item(2, 2).Value = "Efna5"
$ws.Cells.Item(2, 3).Value = "Ephb1"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 1.666083666666667
$ws.Cells.Item(2, 8).Value = 4.998251
$ws.Cells.Item(2, 9).Value = 0.6125276070882968
$ws.Cells.Item(2, 10).Value = 0.6125276070882968
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 0.8299853333333332
$ws.Cells.Item(2, 14).Value = 2.489956
$ws.Cells.Item(2, 15).Value = 0.7350786001848651
$ws.Cells.Item(2, 16).Value = 0.7350786001848651
$ws.Cells.Item(2, 17).Value = 1.382825007439555
$ws.Cells.Item(2, 18).Value = 12.445425066956
$ws.Cells.Item(2, 19).Value = 0.4502559359930503
$ws.Cells.Item(2, 20).Value = 0.4502559359930503

$ws.Cells.Item(3, 1).Value = "FAPs"
$ws.Cells.Item(3, 2).Value = "Efna5"
$ws.Cells.Item(3, 3).Value = "Ephb1"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 1.666083666666667
$ws.Cells.Item(3, 8).Value = 4.998251
$ws.Cells.Item(3, 9).Value = 0.6125276070882968
$ws.Cells.Item(3, 10).Value = 0.6125276070882968
$ws.Cells.Item(3, 11).Value = 2
$ws.Cells.Item(3, 12).Value = 0.6666666666666666
$ws.Cells.Item(3, 13).Value = 0.01102233333333333
$ws.Cells.Item(3, 14).Value = 0.033067
$ws.Cells.Item(3, 15).Value = 0.009761957268446888
$ws.Cells.Item(3, 16).Value = 0.009761957268446888
$ws.Cells.Item(3, 17).Value = 0.01836412953522222
$ws.Cells.Item(3, 18).Value = 0.165277165817
$ws.Cells.Item(3, 19).Value = 0.005979468326139979
$ws.Cells.Item(3, 20).Value = 0.005979468326139979

$ws.Cells.Item(4, 1).Value = "FAPs"
$ws.Cells.Item(4, 2).Value = "Efna5"
$ws.Cells.Item(4, 3).Value = "Ephb1"
$ws.Cells.Item(4, 4).Value = "sCs"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 1.666083666666667
$ws.Cells.Item(4, 8).Value = 4.998251
$ws.Cells.Item(4, 9).Value = 0.6125276070882968
$ws.Cells.Item(4, 10).Value = 0.6125276070882968
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 0.2881033333333333
$ws.Cells.Item(4, 14).Value = 0.8643099999999999
$ws.Cells.Item(4, 15).Value = 0.2551594425466879
$ws.Cells.Item(4, 16).Value = 0.2551594425466879
$ws.Cells.Item(4, 17).Value = 0.4800042579788888
$ws.Cells.Item(4, 18).Value = 4.320038321809999
$ws.Cells.Item(4, 19).Value = 0.1562922027691065
$ws.Cells.Item(4, 20).Value = 0.1562922027691065

$ws.Cells.Item(5, 1).Value = "sCs"
$ws.Cells.Item(5, 2).Value = "Efna5"
$ws.Cells.Item(5, 3).Value = "Ephb1"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 1.053930333333333
$ws.Cells.Item(5, 8).Value = 3.161791
$ws.Cells.Item(5, 9).Value = 0.3874723929117032
$ws.Cells.Item(5, 10).Value = 0.3874723929117031
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 0.8299853333333332
$ws.Cells.Item(5, 14).Value = 2.489956
$ws.Cells.Item(5, 15).Value = 0.7350786001848651
$ws.Cells.Item(5, 16).Value = 0.7350786001848651
$ws.Cells.Item(5, 17).Value = 0.8747467190217777
$ws.Cells.Item(5, 18).Value = 7.872720471196
$ws.Cells.Item(5, 19).Value = 0.2848226641918148
$ws.Cells.Item(5, 20).Value = 0.2848226641918148

$ws.Cells.Item(6, 1).Value = "sCs"
$ws.Cells.Item(6, 2).Value = "Efna5"
$ws.Cells.Item(6, 3).Value = "Ephb1"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 1.053930333333333
$ws.Cells.Item(6, 8).Value = 3.161791
$ws.Cells.Item(6, 9).Value = 0.3874723929117032
$ws.Cells.Item(6, 10).Value = 0.3874723929117031
$ws.Cells.Item(6, 11).Value = 2
$ws.Cells.Item(6, 12).Value = 0.6666666666666666
$ws.Cells.Item(6, 13).Value = 0.01102233333333333
$ws.Cells.Item(6, 14).Value = 0.033067
$ws.Cells.Item(6, 15).Value = 0.009761957268446888
$ws.Cells.Item(6, 16).Value = 0.009761957268446888
$ws.Cells.Item(6, 17).Value = 0.01161677144411111
$ws.Cells.Item(6, 18).Value = 0.104550942997
$ws.Cells.Item(6, 19).Value = 0.00378248894230691
$ws.Cells.Item(6, 20).Value = 0.003782488942306908

$ws.Cells.Item(7, 1).Value = "sCs"
$ws.Cells.Item(7, 2).Value = "Efna5"
$ws.Cells.Item(7, 3).Value = "Ephb1"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 1.053930333333333
$ws.Cells.Item(7, 8).Value = 3.161791
$ws.Cells.Item(7, 9).Value = 0.3874723929117032
$ws.Cells.Item(7, 10).Value = 0.3874723929117031
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 0.2881033333333333
$ws.Cells.Item(7, 14).Value = 0.8643099999999999
$ws.Cells.Item(7, 15).Value = 0.2551594425466879
$ws.Cells.Item(7, 16).Value = 0.2551594425466879
$ws.Cells.Item(7, 17).Value = 0.3036408421344445
$ws.Cells.Item(7, 18).Value = 2.73276757921
$ws.Cells.Item(7, 19).Value = 0.09886723977758143
$ws.Cells.Item(7, 20).Value = 0.09886723977758138

Write-Output "done"